$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates. A leading apostrophe forces text interpretation so
# numeric-looking strings (e.g. "560.64") are not auto-converted to numbers;
# resetting the Style back to Normal afterwards avoids leaving a stray
# quote-prefix / text number-format style on the cell.

$ws.Range("D2").Value = "'66.782.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -5.24%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.376.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -6.39%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.11%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'560.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -5.81%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'184.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -9.24%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -4.26%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.05%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'3.367.14"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -6.30%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -11.76%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.598"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -7.31%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'48.17"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -10.08%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'  -10.21%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'8.76"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -9.65%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.902.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -6.47%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'610.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -11.63%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'66.670.11"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -5.50%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.362.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -6.40%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -4.27%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'17.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -6.79%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -8.22%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.917"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -8.08%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'17.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -7.30%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'5.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.67%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'97.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -11.94%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'4.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -9.93%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -8.97%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'9.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -9.79%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'8.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -11.33%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'30.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -9.86%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "'dogwifhat"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'3.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -13.07%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'NEARProtocol"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'6.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -10.03%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  -8.71%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'Bittensor"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'547.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +7.44%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'Hedera"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'0.106"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -7.18%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'Maker"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'3.841.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.40%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'OKB"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'58.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -7.62%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -0.03%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +38.06%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -5.89%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.0₃0733"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -13.82%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -9.27%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Kaspa"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.129"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -6.17%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'TheGraph"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'0.354"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -7.67%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'32.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -10.69%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0421"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -10.77%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -12.45%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'3.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -9.17%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -6.77%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -0.39%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'7.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -9.94%  "
$ws.Range("E51").Style = "Normal"
